$wb = $excel.ActiveWorkbook

# --- Schedule sheet updates ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 890.0228512499999
$wsSchedule.Range("F2").Value = 19.62131506283069
$wsSchedule.Range("E3").Value = 386.1981435
$wsSchedule.Range("F3").Value = 25.54220525793651

# --- Detailed sheet updates ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B16").Value = 50.62891

$wsDetailed.Range("B17").Value = 36.06
$wsDetailed.Range("C17").Value = "historical"

$wsDetailed.Range("B18").Value = 36.06
$wsDetailed.Range("C18").Value = "historical"

$wsDetailed.Range("B19").Value = 57.06003

$wsDetailed.Range("B21").Value = 36.05971

$wsDetailed.Range("B24").Value = 36.07

$wsDetailed.Range("B27").Value = 36.07
$wsDetailed.Range("B28").Value = 36.07
$wsDetailed.Range("B29").Value = 22.07
$wsDetailed.Range("B30").Value = -5.50985
$wsDetailed.Range("B31").Value = -17.43694
$wsDetailed.Range("B32").Value = -17.36059
$wsDetailed.Range("B33").Value = -15.55074
$wsDetailed.Range("B34").Value = 19.14129
$wsDetailed.Range("B35").Value = 10.3197
$wsDetailed.Range("B36").Value = 2.14574

$wsDetailed.Range("B38").Value = -3.17514
$wsDetailed.Range("B39").Value = 3.14796
$wsDetailed.Range("B40").Value = 0.01138
$wsDetailed.Range("B41").Value = 32.40461
$wsDetailed.Range("B42").Value = 32.40461
$wsDetailed.Range("B43").Value = 32.40461
$wsDetailed.Range("B44").Value = 29.85322
$wsDetailed.Range("B45").Value = 73.20007
$wsDetailed.Range("B46").Value = 57.09

$wsDetailed.Range("B48").Value = 57.0389
$wsDetailed.Range("B49").Value = 57.06003
